# Updates cryptos list values (prices + 1h volume %) for Wed Apr 26 2023 GitHub Actions run.
# Rows 43-46 also swap coin identity (Algorand<->Frax, RenderToken<->WEMIXTOKEN reordered).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: cell address -> new text value. "numeric" cells must be forced to Text
# format first (NumberFormat "@") so Excel stores them as strings, matching the source
# workbook where every Price/Volume cell is an inline string, not a numeric cell.
$updates = [ordered]@{
    "D2" = "29.810.96"
    "E2" = "  +8.58%  "
    "D3" = "1.951.64"
    "E3" = "  +7.01%  "
    "D4" = "1.002"
    "E4" = "  -0.20%  "
    "D5" = "342.42"
    "E5" = "  +2.82%  "
    "E6" = "  -0.23%  "
    "D7" = "0.4769"
    "E7" = "  +4.50%  "
    "D8" = "0.4140"
    "E8" = "  +8.15%  "
    "D9" = "47.93"
    "E9" = "  +4.06%  "
    "D10" = "0.08254"
    "E10" = "  +5.23%  "
    "D11" = "1.036"
    "E11" = "  +8.37%  "
    "D12" = "22.68"
    "E12" = "  +8.02%  "
    "D13" = "1.946.32"
    "E13" = "  +6.87%  "
    "D14" = "6.184"
    "E14" = "  +6.07%  "
    "D15" = "7.417"
    "E15" = "  +5.27%  "
    "D16" = "92.20"
    "E16" = "  +3.14%  "
    "D17" = "1.002"
    "E17" = "  -0.31%  "
    "D18" = "0.00001060"
    "E18" = "  +3.91%  "
    "D19" = "0.06692"
    "E19" = "  +1.60%  "
    "D20" = "18.06"
    "E20" = "  +5.90%  "
    "D21" = "1.001"
    "E21" = "  -0.21%  "
    "D22" = "29.789.87"
    "E22" = "  +8.55%  "
    "D23" = "5.590"
    "E23" = "  +5.94%  "
    "D24" = "11.27"
    "E24" = "  +4.54%  "
    "D25" = "2.262"
    "E25" = "  -0.29%  "
    "D26" = "2.185.76"
    "E26" = "  +7.81%  "
    "D27" = "161.59"
    "E27" = "  +1.45%  "
    "D28" = "20.22"
    "E28" = "  +4.60%  "
    "D29" = "2.187"
    "E29" = "  +7.27%  "
    "D30" = "5.697"
    "E30" = "  +8.29%  "
    "D31" = "122.88"
    "E31" = "  +4.46%  "
    "D32" = "1.013"
    "E32" = "  +9.35%  "
    "D33" = "0.09621"
    "E33" = "  +2.96%  "
    "D34" = "1.477"
    "E34" = "  +12.44%  "
    "D35" = "3.688"
    "E35" = "  +3.27%  "
    "D36" = "5.533"
    "E36" = "  +6.08%  "
    "D37" = "0.06313"
    "E37" = "  +7.05%  "
    "D38" = "0.02323"
    "E38" = "  +6.45%  "
    "D39" = "8.524"
    "E39" = "  +5.59%  "
    "D40" = "1.191"
    "E40" = "  +4.27%  "
    "D41" = "0.6097"
    "E41" = "  +6.71%  "
    "D42" = "10.71"
    "E42" = "  +8.41%  "
    "B43" = "Frax"
    "C43" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
    "D43" = "1.000"
    "E43" = "  -0.17%  "
    "B44" = "Algorand"
    "C44" = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
    "D44" = "0.1890"
    "E44" = "  +4.10%  "
    "B45" = "WEMIXTOKEN"
    "C45" = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
    "D45" = "1.271"
    "E45" = "  -0.37%  "
    "B46" = "RenderToken"
    "C46" = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
    "D46" = "2.365"
    "E46" = "  +30.67%  "
    "D47" = "0.5715"
    "E47" = "  +6.09%  "
    "D48" = "12.53"
    "E48" = "  +6.30%  "
    "D49" = "1.991"
    "E49" = "  +5.92%  "
    "D50" = "0.07341"
    "E50" = "  +5.72%  "
    "D51" = "114.02"
    "E51" = "  +3.52%  "
}

foreach ($addr in $updates.Keys) {
    $value = $updates[$addr]
    $cell = $ws.Range($addr)
    $isNumericLooking = $value -match "^-?\d+(\.\d+)?$"
    if ($isNumericLooking) {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $value
}
